$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.195921421051025
$ws.Range("B1").Value = 2.338304996490479
$ws.Range("C1").Value = 6.780253410339355
$ws.Range("D1").Value = 2.324172735214233
$ws.Range("E1").Value = 1.184273719787598
